# Updated symbol list on Sat Jan 21 14:52:54 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) quotes for the coin rows
# that moved since the last scrape. Values are written as literal text (the
# sheet stores them as plain strings, e.g. "303.21" / "4.96%", not numbers),
# so each cell is forced to Text format before the new quote is written.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

Set-TextValue "D2" '303.21'
Set-TextValue "E2" '4.96%'
Set-TextValue "D3" '34.88'
Set-TextValue "E3" '12.42%'
Set-TextValue "D4" '5.168'
Set-TextValue "E4" '4.28%'
Set-TextValue "D5" '0.07819'
Set-TextValue "E5" '6.31%'
Set-TextValue "D6" '2.339'
Set-TextValue "E6" '2.30%'
Set-TextValue "D7" '8.066'
Set-TextValue "E7" '4.41%'
Set-TextValue "D8" '3.979'
Set-TextValue "E8" '6.57%'
Set-TextValue "D9" '0.9258'
Set-TextValue "E9" '1.56%'
Set-TextValue "D10" '0.1010'
Set-TextValue "E10" '8.56%'
Set-TextValue "D11" '0.1825'
Set-TextValue "E11" '7.07%'
Set-TextValue "D12" '0.08540'
Set-TextValue "E12" '3.39%'
Set-TextValue "D13" '0.03418'
Set-TextValue "E13" '10.27%'
Set-TextValue "D14" '0.09904'
Set-TextValue "E14" '-0.73%'
Set-TextValue "D15" '0.001500'
Set-TextValue "E15" '0.39%'
Set-TextValue "D16" '0.005730'
Set-TextValue "E16" '-0.26%'
Set-TextValue "E17" '0.16%'
Set-TextValue "D18" '2.109'
Set-TextValue "E18" '3.29%'
Set-TextValue "E19" '2.88%'
Set-TextValue "D20" '0.1325'
Set-TextValue "E20" '2.85%'
Set-TextValue "D21" '4.566'
Set-TextValue "E21" '9.31%'
Set-TextValue "D23" '0.04645'
Set-TextValue "E23" '2.76%'
Set-TextValue "D24" '0.001217'
Set-TextValue "E24" '0.49%'
Set-TextValue "D25" '0.004340'
Set-TextValue "E25" '3.68%'
Set-TextValue "D26" '0.0001302'
Set-TextValue "D27" '0.0003395'
Set-TextValue "E27" '0.12%'
Set-TextValue "E39" '11.69%'
Set-TextValue "D40" '0.04748'
Set-TextValue "E40" '5.95%'
Set-TextValue "D41" '0.007759'
Set-TextValue "E41" '5.03%'
Set-TextValue "D42" '0.1410'
Set-TextValue "E42" '5.80%'
Set-TextValue "D43" '0.008512'
Set-TextValue "E43" '-13.45%'
Set-TextValue "D44" '0.002293'
Set-TextValue "E44" '7.22%'
Set-TextValue "D45" '0.01003'
Set-TextValue "E45" '11.66%'
Set-TextValue "D46" '0.00006087'
Set-TextValue "E46" '-0.12%'
Set-TextValue "D47" '0.00000000749'
Set-TextValue "E47" '-0.06%'
Set-TextValue "D48" '3.909'
Set-TextValue "E48" '52.37%'
Set-TextValue "D49" '0.002687'
Set-TextValue "E49" '28.11%'
Set-TextValue "D50" '0.00002098'
Set-TextValue "E50" '-0.06%'
Set-TextValue "D51" '0.0001998'
Set-TextValue "E51" '-0.06%'
